$wb = $excel.ActiveWorkbook

$wsSettings  = $wb.Worksheets.Item("Settings")
$wsConstants = $wb.Worksheets.Item("Constants")

# ---------------------------------------------------------------------------
# Settings sheet: update existing values, then add the new System1 / Exception
# rows (order below reproduces the exact shared-string insertion order
# recorded in the target workbook).
# ---------------------------------------------------------------------------
$wsSettings.Range("B2").Value = "GenerateYearlyReports"
$wsSettings.Range("B5").Value = "GenerateYearlyReport-Dispatcher"

$wsSettings.Range("B6").Value = "https://acme-test.uipath.com/work-items"
$wsSettings.Range("A7").Value = "System1_URL"
$wsSettings.Range("B7").Value = "https://acme-test.uipath.com"
$wsSettings.Range("A8").Value = "System1_CredentialName"
$wsSettings.Range("B8").Value = "ACMELogin"
$wsSettings.Range("A9").Value = "ExceptionEmail"
$wsSettings.Range("B9").Value = "exceptions@acme-test.com"
$wsSettings.Range("A10").Value = "System1_WorkItemType"
$wsSettings.Range("B10").Value = "WI4"
$wsSettings.Range("A6").Value = "System1_WorkItemsURL"

# Hyperlinks for the URL / e-mail cells (Excel auto-applies the "Hyperlink"
# cell style to these cells).
$wsSettings.Hyperlinks.Add($wsSettings.Range("B6"), "https://acme-test.uipath.com/work-items")
$wsSettings.Hyperlinks.Add($wsSettings.Range("B7"), "https://acme-test.uipath.com")
$wsSettings.Hyperlinks.Add($wsSettings.Range("B9"), "mailto:exceptions@acme-test.com")

# ---------------------------------------------------------------------------
# Constants sheet: MaxRetryNumber 0 -> 2
# ---------------------------------------------------------------------------
$wsConstants.Range("B2").Value = 2

# ---------------------------------------------------------------------------
# Restore the cursor / selection state recorded in the target workbook.
# ---------------------------------------------------------------------------
$wsConstants.Activate()
$wsConstants.Range("B3").Select()

$wsSettings.Activate()
$wsSettings.Range("A7").Select()

Write-Host "Config.xlsx updated"
